$wb = $excel.ActiveWorkbook

# --- BD sheet: add a second bridge domain row ---
$bd = $wb.Worksheets.Item("BD")
$bd.Range("A3").Value = "bd"
$bd.Range("B3").Value = "mark_second_bd_for_subnet"
$bd.Range("C3").Value = "This other  bridge domain is created by the Terraform ACI provider1"
$bd.Columns.Item(3).ColumnWidth = 56.09765625
$bd.Range("B3").Select()

# --- EPG sheet: add a second EPG row referencing the new BD ---
$epg = $wb.Worksheets.Item("EPG")
$epg.Range("A3").Value = "epg"
$epg.Range("B3").Value = "mark_second_epg"
$epg.Range("C3").Value = "mark_second_bd_for_subnet"
$epg.Range("D3").Value = "3tier_app"

# Make EPG the active (selected) sheet/tab, with D3 selected
$epg.Activate()
$epg.Range("D3").Select()
